$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target state of the table body (A2:C207): rows re-grouped/re-ordered
# by katedra+zkratka while every (nazev, zkratka, cviciciUcitIdno) triple from
# the original sheet is preserved exactly (same multiset, just reshuffled).
$data = @(
  ,("Fyzika pevných látek I","K839",69)
  ,("Základy počítačové fyziky","K503",69)
  ,("Demografie","5DEMO",130)
  ,("Demografie","FDEMO",130)
  ,("Geometrie","P335",269)
  ,("Geometrie","K312",269)
  ,("Fyziologie živočichů a člověka","PB502",281)
  ,("Diplomová práce I","M200",281)
  ,("Diplomová práce I","N028",281)
  ,("Fyziologie a mol. bio. eukar. buňky I","N016",281)
  ,("Fyziologie živočichů a člověka","KB502",281)
  ,("Průběžná pedagogická praxe ZŠ","P732",302)
  ,("Exkurze","P841",302)
  ,("Proseminář řešení fyzikálních úloh","K622",302)
  ,("Fyzikální praktikum C","K521",302)
  ,("Didaktická a laboratorní technika","PA31",302)
  ,("Proseminář řešení fyzikálních úloh A ZŠ","P736",302)
  ,("Proseminář řešení fyzikálních úloh A","P622",302)
  ,("Souvislá pedagogická praxe II - ZŠ","KA35",302)
  ,("Matematický proseminář","P365",306)
  ,("Diplomový seminář","K505",306)
  ,("Vybrané partie z fyziky","K414",306)
  ,("Teorie elektromagnetického pole","P508",307)
  ,("Seminář katedry","P602",307)
  ,("Seminář katedry","K602",307)
  ,("PVK - Vyb.partie z teoretické fyziky I","K422",307)
  ,("Souhrnný seminář - počítačové modelování","KSPM",307)
  ,("Vybrané partie z teoretické fyziky I","P938",307)
  ,("Kvantová mechanika","KKM",307)
  ,("Teorie elektromagnetického pole","TEP",307)
  ,("Seminář katedry I","P952",307)
  ,("Odborný seminář katedry A","P531",307)
  ,("Kvantová fyzika","P737",307)
  ,("Kvantová mechanika","KM",307)
  ,("Vybrané partie z fyziky","P427",307)
  ,("Vybrané partie z teorie elektromag. pole","PD203",307)
  ,("PVK - Vyb.partie z teoretické fyziky II","K522",307)
  ,("Fyzika I","K222",307)
  ,("Kvantová fyzika I","K512",307)
  ,("Kvantová fyzika","K724",307)
  ,("Souhrnný seminář - počítačové modelování","SPM",307)
  ,("Kvantová fyzika I","P509",307)
  ,("Kvantová fyzika II","P607",307)
  ,("Seminář katedry II","PA52",307)
  ,("Potenciál kulturní krajiny a GIS","M214",313)
  ,("Metody geografického výzkumu krajiny","M409",313)
  ,("VK - Geografická exkurze v ČR (5 dní)","0101",317)
  ,("Dějiny chemie","KN16",336)
  ,("Didaktika chemie II - SŠ","KN22",336)
  ,("Didaktika chemie II - SŠ","N304",336)
  ,("Experimentální využití elmg. vln B","0218",605)
  ,("Řešení obvodu RLC","0214",605)
  ,("Experimentální využití elmg. vln F","0237",605)
  ,("Vybrané matematické metody ve fyzice","0212",605)
  ,("VK- Praktikum mikrovlnné techniky D","0211",605)
  ,("Praktikum mikrovlnné techniky - aplikace","0207",605)
  ,("VK-Programování ATMEL I","0175",605)
  ,("Teoretická mechanika I","K305",606)
  ,("Vybrané partie z teoretické fyziky II","P522",606)
  ,("Vybrané partie z teoretické fyziky I","P422",606)
  ,("Počítačové modelování III","P401",612)
  ,("PVK - Hromadné zpracování dat","K406",612)
  ,("Programování A","K103",612)
  ,("Programování B","P203",612)
  ,("PVK - Hromadné zpracování dat","P406",612)
  ,("Počítačové modelování I","K107",612)
  ,("VK-Fyzikální seminář","0204",612)
  ,("Programování A","P103",612)
  ,("Počítačové modelování II","K206",612)
  ,("Počítačové modelování II","P206",612)
  ,("Systém kvality a vyhodnocování dat","N111",719)
  ,("Počítač. modelování - spojité modelování","P366",816)
  ,("VK - Exkurze z regionální geologie","0109",827)
  ,("Základy financí","ZFIN",929)
  ,("Finanční analýza","FINA",929)
  ,("Základy matematiky","K103",1108)
  ,("Finanční analýza","FINA",1383)
  ,("Statistika na PC","K516",1456)
  ,("Matematika IV","K409",1537)
  ,("Somatologie","BP423",1546)
  ,("Obecná parazitologie","P323",1546)
  ,("Biologie a ekologie člověka I","P303",1546)
  ,("Obecná zoologie","K110",1546)
  ,("Biologie parazitů","N032",1546)
  ,("Somatologie","BK423",1546)
  ,("Ekofyziologie","N051",1546)
  ,("Komplexní analýza materiálů","PD202",1590)
  ,("Kvantová fyzika I","P509",1590)
  ,("Met. analýz mat. `nvyuž. el. mikroskopie","AP03",1590)
  ,("Metody analýz materiálů","AP05",1590)
  ,("Metody molekulární biologie","N023",1986)
  ,("Metody molekulární biologie","ME200",1986)
  ,("Algologie","P420",1999)
  ,("Limnobiologie","P328",1999)
  ,("Environmentální legislativa","K507",1999)
  ,("Limnobiologie","P528",1999)
  ,("Limnobiologie","K311",1999)
  ,("Limnobiologie","BK108",1999)
  ,("Environmentální legislativa","P525",1999)
  ,("Algologie","N036",1999)
  ,("Didaktika geografie I","N001",2232)
  ,("Didaktika geografie I","M105",2232)
  ,("Labor. Techniq. in Mol. and Cell Biology","E121",2441)
  ,("Microscopy Techniques in Biology","E120",2441)
  ,("Eliminace polutantů z fluidních systémů","N021",2450)
  ,("Anglický jazyk II","K411",2591)
  ,("Biologie a ekologie člověka I","P303",2757)
  ,("Souvislá pedagogická praxe I SŠ","P706",2757)
  ,("VK-Česká republika a EU","0116",2758)
  ,("Geospatial data: sources and analysis","PD107",2787)
  ,("VK-Němčina pro geografy","0126",2812)
  ,("VK - Konverzace němčiny pro geografy","0130",2812)
  ,("Metody molekulární biologie","N023",2997)
  ,("Metody molekulární biologie","ME200",2997)
  ,("Medicína katastrof","P518",3066)
  ,("Technologie ochrany životního prostředí","P414",3175)
  ,("Toxikologie rostl. a živočišných toxinů","P534",3176)
  ,("Ekotoxikologie","N109",3176)
  ,("Počítačová kriminalita","0141",3458)
  ,("Welfare Assessment of zoo animals","E133",3796)
  ,("Preparace a konzervování přírodnin II","K403",3797)
  ,("Preparace a konzervování přírodnin II","P416",3797)
  ,("Preparace a konzervování přírodnin II","P419",3797)
  ,("Preparace a konzervování přírodnin II","K403",3798)
  ,("Preparace a konzervování přírodnin II","P416",3798)
  ,("Preparace a konzervování přírodnin II","P419",3798)
  ,("Funkce a struktura buněčných membrán","0088",4175)
  ,("Základy receptorové neurofyziologie","0089",4176)
  ,("Základy bioinženýrství","0090",4177)
  ,("Modelování podnikových procesů","0160",4190)
  ,("Počítačové sítě a protokoly","PSP",4195)
  ,("Fyziologie a mol. bio. eukar. buňky I","N016",4208)
  ,("Matematika I","K106",4221)
  ,("Úvod do fyziky","P112",4221)
  ,("Matematika I","P106",4221)
  ,("Kvantová fyzika I","P509",4221)
  ,("Neziskový sektor a regionální rozvoj","0150",4227)
  ,("Moderní elektroanal. met. - prakt. kurz","N029",4258)
  ,("PVK - Fyzikálně chemická cvičení A","P506",4258)
  ,("Fyzikálně chemická cvičení","P517",4258)
  ,("Laboratorní cvičení z fyzikální chemie","P532",4258)
  ,("Fyzikálně chemická cvičení","K502",4258)
  ,("Fyzikálně chemická cvičení","P527",4258)
  ,("Vybrané kapitoly z geografie krajiny","0152",4302)
  ,("Anglický jazyk A","K100",4380)
  ,("Explantátové kultury rostlin I","N049",4390)
  ,("Právní aspekty vzdělávání a předpisy","KN46",4581)
  ,("Právní aspekty vzdělávání a předpisy","N005",4581)
  ,("Information Security","0177",4685)
  ,("Mikroprocesory a senzory v praxi I","0164",4746)
  ,("Mikroprocesory a senzory v praxi II","0171",4746)
  ,("Scientific inquiry and reasoning","PD101",4747)
  ,("Fyzikálně chemická cvičení","P517",4757)
  ,("Automatické řízení","AUC",5039)
  ,("Chemické a biologické ohrožení","N020",5062)
  ,("Ochrana přírody","P312",5131)
  ,("Praktikum validace analytických metod","M208",5135)
  ,("Fyzika pevných látek I","P861",5201)
  ,("Fyzika pevných látek I","P839",5201)
  ,("Lokální a regionální rozvoj","M101",5240)
  ,("Metody geografického výzkumu krajiny","M203",5240)
  ,("Životní prostředí a udržitelný rozvoj","BK401",5240)
  ,("English for Science","E403",5240)
  ,("Životní prostředí a udržitelný rozvoj","B401",5240)
  ,("Terénní cvičení ze SG","M312",5240)
  ,("Zoologické terénní cvičení","KB427",5472)
  ,("Zoologické terénní cvičení","PB427",5472)
  ,("Ornitologie","MA205",5472)
  ,("Srovnávací anatom. a morfol. obratlovců","N060",5472)
  ,("Anglický jazyk A","K100",5754)
  ,("Biomedicínské mat.- přípr., char., a vl.","P328",5854)
  ,("Identif. a hodnocení ekosystém. služeb","M317",5886)
  ,("Geografie a dědictví krajiny","0156",5924)
  ,("Environmental Biotechnology","M303",5927)
  ,("Základy financí","ZFIN",6120)
  ,("Nekovové anorganické materiály","AMAT",6886)
  ,("Nekovové anorganické materiály","AMAT",6886)
  ,("Úvod do modelování průmyslových procesů","M211",7307)
  ,("Odběr, příprava a analýza vzorků","N112",7416)
  ,("Přípr. a využití nanovlákenných struktur","PD55",7805)
  ,("Didaktika biologie II","P208",8021)
  ,("Didaktika biologie I","P108",8021)
  ,("Doučování","DOUC",8021)
  ,("Souvislá pedagogická praxe I - SŠ","P707A",8021)
  ,("Průběžná pedagogická praxe","P705",8021)
  ,("Souvislá pedagogická praxe II - SŠ","P707B",8021)
  ,("Kompozity nanočástic kovů a jejich oxidů","PD12",8074)
  ,("Nanokompozitní materiály na bázi uhlíku","PD08",8075)
  ,("Apl. biologie prokaryot. a eukaryot. m.","ME301",8154)
  ,("Biologie půdních mikroorganismů","M101",8154)
  ,("Kvantová mechanika","KKM",8310)
  ,("Kvantová mechanika","KM",8310)
  ,("Microscopy Techniques in Biology","E120",8316)
  ,("Algoritmizace a programování I","APR1",8514)
  ,("Algoritmizace a programování II","APR2",8514)
  ,("Machine Learning Based on R Software","EMLR",8514)
  ,("Data Analysis and Visualisation","EDAV",8514)
  ,("Python and R for Data Science","EPYR",8514)
  ,("Data Mining Techniq. based on R Software","EDMR",8514)
  ,("Základy matematiky","P103",8677)
  ,("Pravděpodobnost a statistika I","K413",8677)
  ,("Vybrané partie z matematiky","P232",8677)
  ,("Podniková ekonomika","PEK",8749)
  ,("Pravděpodobnost a statistika II","K512",8750)
  ,("Objektově orientované návrhové vzory","OONV",8753)
  ,("Objektově orientované návrhové vzory","KOONV",8753)
)

$nRows = 206
$arr = New-Object 'object[,]' $nRows,3
for ($i = 0; $i -lt $nRows; $i++) {
  $r = $data[$i]
  $arr[$i,0] = $r[0]
  $arr[$i,1] = $r[1]
  $arr[$i,2] = $r[2]
}
$lastRow = 1 + $nRows
$rng = $ws.Range("A2:C" + $lastRow)
$rng.Value2 = $arr

# "zkratka" values that look numeric (e.g. "0101") get mangled by the generic
# array write above (Excel auto-converts them to a number and drops the leading
# zero). Re-enter those specific cells as Text so they keep their original string form.
$textCells = @(48,52,53,54,55,56,57,58,68,74,110,112,113,120,128,129,130,131,138,145,150,151,152,174)
foreach ($r in $textCells) {
  $cell = $ws.Cells.Item($r, 2)
  $cell.NumberFormat = "@"
  $cell.Value2 = $arr[$r - 2, 1]
}

Write-Host "Wrote $nRows rows"
